$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column F (dSF) for rows 2-26, per the diff.
$newValues = @{
    2  = -1
    3  = -1
    5  = -2
    6  = 2
    7  = -2
    8  = 7
    9  = 1
    11 = 2
    12 = 1
    15 = 6
    16 = -1
    17 = 9
    18 = -3
    19 = 1
    20 = -2
    22 = -2
    23 = 1
    24 = -2
    25 = 1
    26 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $newValues[$row]
}
